$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new "Wins" / "Losses" / "Ties" columns ---------
# Clone the existing header formatting (bold, centered, bordered) from the
# last header cell (AC1) onto the three new header cells, then set text.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2-48): season record for every player on the roster -----
$lastRow = 48
$wins = $ws.Range("AD2:AD$lastRow")
$losses = $ws.Range("AE2:AE$lastRow")
$ties = $ws.Range("AF2:AF$lastRow")

$wins.Value = 67
$losses.Value = 95
$ties.Value = 0
